$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 63, pushing existing rows 63..131 down to 64..132
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly record
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 44902
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102004
$ws.Range("J63").Value = "Mandarina"
$ws.Range("K63").Value = "Murcott"
$ws.Range("L63").Value = "Segunda"
$ws.Range("M63").Value = 300
$ws.Range("N63").Value = 17000
$ws.Range("O63").Value = 18000
$ws.Range("P63").Value = 17500
$ws.Range("Q63").Value = "$/caja 20 kilos"
$ws.Range("R63").Value = "Región de Coquimbo"
$ws.Range("S63").Value = 875
$ws.Range("T63").Value = 20
